$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the product_item names in column A (shared strings "tested 13/14/15" -> "tested 16/17/188")
$ws.Range("A2").Value = "tested 16"
$ws.Range("A3").Value = "tested 17"
$ws.Range("A4").Value = "tested 188"

# Move the active selection from E13 to D15 (navigation to accounting module)
$ws.Range("D15").Select()
